# Auto-generated edit script applying Seraph_Profits.xlsx market-price update
# (scheduled runner refresh of currentAveragePrice / LevePrice / LeveProfit columns)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3997.1428
$ws.Range("I74").Value = 3997.1428
$ws.Range("K74").Value = 3997.1428
$ws.Range("M74").Value = -3061.1428
$ws.Range("H77").Value = 3997.1428
$ws.Range("I77").Value = 3997.1428
$ws.Range("K77").Value = 19985.714
$ws.Range("M77").Value = -15305.714
$ws.Range("H113").Value = 9559.6
$ws.Range("I113").Value = 9266.333000000001
$ws.Range("K113").Value = 9266.333000000001
$ws.Range("M113").Value = -6012.333000000001
$ws.Range("H138").Value = 3230.4324
$ws.Range("I138").Value = 2070.25
$ws.Range("J138").Value = 3550.4827
$ws.Range("K138").Value = 6210.75
$ws.Range("L138").Value = 10651.4481
$ws.Range("M138").Value = -1070.75
$ws.Range("N138").Value = -20931.4481

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9557.223
$ws.Range("I61").Value = 10064.375
$ws.Range("J61").Value = 5500
$ws.Range("K61").Value = 10064.375
$ws.Range("L61").Value = 5500
$ws.Range("M61").Value = -9852.375
$ws.Range("N61").Value = -5924
$ws.Range("H76").Value = 35999.5
$ws.Range("J76").Value = 35999.5
$ws.Range("L76").Value = 35999.5
$ws.Range("N76").Value = -36675.5
$ws.Range("H79").Value = 35999.5
$ws.Range("J79").Value = 35999.5
$ws.Range("L79").Value = 35999.5
$ws.Range("N79").Value = -38339.5
$ws.Range("H132").Value = 9317.166999999999
$ws.Range("I132").Value = 10380.6
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 31141.8
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -28611.8
$ws.Range("N132").Value = -17060
$ws.Range("H136").Value = 9557.223
$ws.Range("I136").Value = 10064.375
$ws.Range("J136").Value = 5500
$ws.Range("K136").Value = 30193.125
$ws.Range("L136").Value = 16500
$ws.Range("M136").Value = -27643.125
$ws.Range("N136").Value = -21600

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 13825
$ws.Range("J88").Value = 13825
$ws.Range("L88").Value = 13825
$ws.Range("N88").Value = -14637
$ws.Range("H91").Value = 13825
$ws.Range("J91").Value = 13825
$ws.Range("L91").Value = 13825
$ws.Range("N91").Value = -16633

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3448.2812
$ws.Range("I31").Value = 2387.8262
$ws.Range("J31").Value = 6158.3335
$ws.Range("K31").Value = 2387.8262
$ws.Range("L31").Value = 6158.3335
$ws.Range("M31").Value = -2092.8262
$ws.Range("N31").Value = -6748.3335
$ws.Range("H34").Value = 3448.2812
$ws.Range("I34").Value = 2387.8262
$ws.Range("J34").Value = 6158.3335
$ws.Range("K34").Value = 2387.8262
$ws.Range("L34").Value = 6158.3335
$ws.Range("M34").Value = -2185.8262
$ws.Range("N34").Value = -6562.3335
$ws.Range("H132").Value = 2224.0833
$ws.Range("I132").Value = 2068.9
$ws.Range("K132").Value = 6206.700000000001
$ws.Range("M132").Value = -3676.700000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1725
$ws.Range("I81").Value = 300
$ws.Range("K81").Value = 900
$ws.Range("M81").Value = 223
$ws.Range("H84").Value = 1725
$ws.Range("I84").Value = 300
$ws.Range("K84").Value = 2700
$ws.Range("M84").Value = 2916
$ws.Range("H132").Value = 3487.9
$ws.Range("I132").Value = 3140.5715
$ws.Range("J132").Value = 4298.3335
$ws.Range("K132").Value = 28265.1435
$ws.Range("L132").Value = 38685.0015
$ws.Range("M132").Value = -25735.1435
$ws.Range("N132").Value = -43745.0015
$ws.Range("H137").Value = 6250
$ws.Range("I137").Value = 7500
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 22500
$ws.Range("L137").Value = 15000
$ws.Range("M137").Value = -17400
$ws.Range("N137").Value = -25200

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 10500000
$ws.Range("I3").Value = 10500000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 10500000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -10499884
$ws.Range("N3").ClearContents()
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("H102").Value = 1973.1428
$ws.Range("J102").Value = 2000
$ws.Range("L102").Value = 2000
$ws.Range("N102").Value = -5244
$ws.Range("H132").Value = 4997.8
$ws.Range("I132").Value = 4997.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 14993.4
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12463.4
$ws.Range("N132").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 53666.332
$ws.Range("I18").Value = 47999.5
$ws.Range("J18").Value = 65000
$ws.Range("K18").Value = 47999.5
$ws.Range("L18").Value = 65000
$ws.Range("M18").Value = -47827.5
$ws.Range("N18").Value = -65344
$ws.Range("H22").Value = 17999.908
$ws.Range("I22").Value = 19333.334
$ws.Range("K22").Value = 19333.334
$ws.Range("M22").Value = -19038.334
$ws.Range("H27").Value = 17999.908
$ws.Range("I27").Value = 19333.334
$ws.Range("K27").Value = 19333.334
$ws.Range("M27").Value = -19226.334
$ws.Range("H40").Value = 4568.4
$ws.Range("I40").Value = 4473
$ws.Range("K40").Value = 4473
$ws.Range("M40").Value = -4337
$ws.Range("H60").Value = 50000
$ws.Range("J60").Value = 50000
$ws.Range("L60").Value = 50000
$ws.Range("N60").Value = -51018
$ws.Range("H122").Value = 6957.3335
$ws.Range("I122").Value = 6957.3335
$ws.Range("K122").Value = 20872.0005
$ws.Range("M122").Value = -18422.0005
$ws.Range("H132").Value = 127192.375
$ws.Range("I132").Value = 145045.28
$ws.Range("K132").Value = 435135.84
$ws.Range("M132").Value = -432605.84
$ws.Range("H136").Value = 5504.4116
$ws.Range("J136").Value = 6844.8335
$ws.Range("L136").Value = 20534.5005
$ws.Range("N136").Value = -25634.5005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 700
$ws.Range("I17").Value = 450
$ws.Range("K17").Value = 450
$ws.Range("M17").Value = -278
$ws.Range("H81").Value = 6594.6
$ws.Range("I81").Value = 5135.2856
$ws.Range("K81").Value = 10270.5712
$ws.Range("M81").Value = -9209.5712
$ws.Range("H84").Value = 6594.6
$ws.Range("I84").Value = 5135.2856
$ws.Range("K84").Value = 51352.856
$ws.Range("M84").Value = -46048.856
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H132").Value = 1735.6666
$ws.Range("I132").Value = 1677.2222
$ws.Range("K132").Value = 5031.6666
$ws.Range("M132").Value = -2501.6666
$ws.Range("H136").Value = 2388.3333
$ws.Range("I136").Value = 1298.1111
$ws.Range("J136").Value = 4023.6667
$ws.Range("K136").Value = 3894.3333
$ws.Range("L136").Value = 12071.0001
$ws.Range("M136").Value = -1344.3333
$ws.Range("N136").Value = -17171.0001
